$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I1:J1 header cells should carry the same style as the other header cells (s="1").
# Copy H1's format onto I1:J1, then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns: I = I0, J = IF, for rows 2-78.
$ijValues = @{
    2 = @(4, 5)
    3 = @(6, 6)
    4 = @(3, 5)
    5 = @(5, 6)
    6 = @(6, 7)
    7 = @(6, 6)
    8 = @(5, 5)
    9 = @(6, 7)
    10 = @(6, 7)
    11 = @(7, 7)
    12 = @(5, 5)
    13 = @(5, 5)
    14 = @(9, 9)
    15 = @(7, 7)
    16 = @(6, 6)
    17 = @(10, 10)
    18 = @(6, 6)
    19 = @(5, 6)
    20 = @(6, 6)
    21 = @(6, 7)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(6, 6)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(7, 8)
    29 = @(7, 7)
    30 = @(7, 7)
    31 = @(6, 6)
    32 = @(7, 7)
    33 = @(8, 8)
    34 = @(6, 6)
    35 = @(7, 7)
    36 = @(7, 7)
    37 = @(7, 7)
    38 = @(8, 8)
    39 = @(8, 8)
    40 = @(6, 6)
    41 = @(7, 7)
    42 = @(7, 7)
    43 = @(8, 8)
    44 = @(8, 8)
    45 = @(7, 7)
    46 = @(6, 6)
    47 = @(8, 8)
    48 = @(7, 7)
    49 = @(7, 7)
    50 = @(7, 7)
    51 = @(7, 7)
    52 = @(8, 8)
    53 = @(8, 8)
    54 = @(8, 8)
    55 = @(11, 11)
    56 = @(8, 8)
    57 = @(9, 9)
    58 = @(8, 8)
    59 = @(8, 8)
    60 = @(8, 8)
    61 = @(7, 8)
    62 = @(7, 7)
    63 = @(8, 8)
    64 = @(9, 9)
    65 = @(8, 8)
    66 = @(8, 8)
    67 = @(6, 6)
    68 = @(8, 8)
    69 = @(7, 7)
    70 = @(7, 7)
    71 = @(7, 7)
    72 = @(8, 8)
    73 = @(8, 8)
    74 = @(5, 5)
    75 = @(5, 5)
    76 = @(8, 8)
    77 = @(7, 7)
    78 = @(5, 5)
}

foreach ($rowNum in $ijValues.Keys) {
    $pair = $ijValues[$rowNum]
    $ws.Cells.Item([int]$rowNum, 9).Value = $pair[0]
    $ws.Cells.Item([int]$rowNum, 10).Value = $pair[1]
}
